$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '37.895.73'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.08%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.035.96'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.68%  '

$ws.Range("E4").Value = '  -0.05%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '227.55'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("E6").Value = '  -0.57%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '60.38'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +3.10%  '

$ws.Range("E8").Value = '  -0.01%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.382'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '

$ws.Range("E10").Value = '  +1.15%  '

$ws.Range("E11").Value = '  +0.79%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '2.336.91'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.72%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '14.57'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.58%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.33'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.58%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.761'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.02%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '5.15'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.63%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.033.20'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.58%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '37.872.34'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.01%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '69.87'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.23%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.93'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -5.27%  '

$ws.Range("E21").Value = '  -1.12%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '224.75'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("E23").Value = '  +0.10%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.80%  '

$ws.Range("E25").Value = '  +1.05%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '166.91'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.33%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.30'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("E28").Value = '  -3.51%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '18.92'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.49%  '

$ws.Range("E30").Value = '  -4.21%  '

$ws.Range("E31").Value = '  +1.29%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.14'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +4.72%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.42'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.35%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.00%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0606'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.84%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.44'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +5.75%  '

$ws.Range("E37").Value = '  -2.44%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.26'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.85%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.06%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.525.15'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.89%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '17.11'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +3.72%  '

$ws.Range("E42").Value = '  +0.45%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '96.14'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.00%  '

$ws.Range("E44").Value = '  -0.42%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0914'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.82%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.11'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.42%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.00'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.67%  '

$ws.Range("E48").Value = '  -0.50%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("E50").Value = '  +0.48%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.224.65'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.76%  '
